$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 4) renames ---
$ws.Range("E4").Value = "Sub_Name"
$ws.Range("F4").Value = "Days_Left"

# --- Data rows 5-9: Sub_Name (col C) and Sub_ID (col D) rework ---
$ws.Range("C5").Value = "Steven"
$ws.Range("D5").Value = "00002"
$ws.Range("E5").Value = "John Doe"

$ws.Range("C6").Value = "John Doe"
$ws.Range("D6").Value = "00004"
$ws.Range("E6").Value = "Max Mad"

$ws.Range("C7").Value = "Jack D Ribba"
$ws.Range("D7").Value = "00005"
$ws.Range("E7").Value = "Test"

$ws.Range("C8").Value = "Max Mad"
$ws.Range("D8").Value = "00002"
$ws.Range("E8").Value = "John Doe"

$ws.Range("C9").Value = "Test"
$ws.Range("D9").Value = "00003"
$ws.Range("E9").Value = "Jack D Ribba"

# --- Remove the now-unused trailing blank formatted row ---
$ws.Rows("30:30").Delete()

# --- View / window state ---
$excel.ActiveWindow.Zoom = 145

$sel = $ws.Range("B12")
[void]$sel.Select()
